# Adjusted excel for start-up/shut-down limits.
#
# 1) A new unit__to_node relationship row (Solar_Plant_Kasso -> Power_Kasso)
#    is inserted into "Object__to_from_node" for start_up_limit (0.5) and
#    shut_down_limit (0.8), right after the existing unit_capacity row for
#    that same relationship. All rows below shift down by two.
# 2) As a consequence of that edit (re-saved from the source notebook), the
#    "Nodes" table - and the matching "Definition" object list - ended up
#    re-ordered (same node set/attributes, different row order).
#
# NOTE: reading `.Value` back out of a Range in this host returns the
# reflected property stub instead of the cell contents, so all reads below
# go through `.Value2` (writes use `.Value2` too, for symmetry).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Object__to_from_node: insert the two new start-up/shut-down rows.
# ---------------------------------------------------------------------
$objSheet = $wb.Worksheets.Item("Object__to_from_node")

# Snapshot rows 3..20 first (old layout) before writing anything, then
# replay them two rows down (new layout rows 5..22).
$savedRows = @()
for ($r = 3; $r -le 20; $r++) {
    $savedRows += ,@(
        $objSheet.Range("A$r").Value2,
        $objSheet.Range("B$r").Value2,
        $objSheet.Range("C$r").Value2,
        $objSheet.Range("D$r").Value2,
        $objSheet.Range("E$r").Value2,
        $objSheet.Range("F$r").Value2
    )
}

for ($i = 0; $i -lt $savedRows.Length; $i++) {
    $dest = $i + 5
    $row = $savedRows[$i]
    $objSheet.Range("A$dest").Value2 = $row[0]
    $objSheet.Range("B$dest").Value2 = $row[1]
    $objSheet.Range("C$dest").Value2 = $row[2]
    $objSheet.Range("D$dest").Value2 = $row[3]
    $objSheet.Range("E$dest").Value2 = $row[4]
    $objSheet.Range("F$dest").Value2 = $row[5]
}

# New row 3: Solar_Plant_Kasso -> Power_Kasso, start_up_limit = 0.5
$objSheet.Range("A3").Value2 = "unit__to_node"
$objSheet.Range("B3").Value2 = "unit"
$objSheet.Range("C3").Value2 = "Solar_Plant_Kasso"
$objSheet.Range("D3").Value2 = "Power_Kasso"
$objSheet.Range("E3").Value2 = "start_up_limit"
$objSheet.Range("F3").Value2 = 0.5

# New row 4: Solar_Plant_Kasso -> Power_Kasso, shut_down_limit = 0.8
$objSheet.Range("A4").Value2 = "unit__to_node"
$objSheet.Range("B4").Value2 = "unit"
$objSheet.Range("C4").Value2 = "Solar_Plant_Kasso"
$objSheet.Range("D4").Value2 = "Power_Kasso"
$objSheet.Range("E4").Value2 = "shut_down_limit"
$objSheet.Range("F4").Value2 = 0.8

$objSheet.Range("A1:F22").EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 2) Nodes: re-order the 12 node rows (same data, new row order).
# ---------------------------------------------------------------------
$nodesSheet = $wb.Worksheets.Item("Nodes")

$nodesOrder = @(
    "Power_Kasso",
    "Hydrogen_storage_Kasso",
    "District_Heating",
    "Raw_Methanol",
    "Carbon_Dioxide",
    "Hydrogen_Kasso",
    "Waste_Heat",
    "Water",
    "Vaporized_Carbon_Dioxide",
    "E-Methanol_Kasso",
    "E-Methanol_storage_Kasso",
    "Power_Wholesale"
)

$nodeData = @{}
for ($r = 2; $r -le 13; $r++) {
    $name = $nodesSheet.Range("A$r").Value2
    $nodeData[$name] = @(
        $nodesSheet.Range("B$r").Value2,
        $nodesSheet.Range("C$r").Value2,
        $nodesSheet.Range("D$r").Value2,
        $nodesSheet.Range("E$r").Value2,
        $nodesSheet.Range("F$r").Value2,
        $nodesSheet.Range("G$r").Value2
    )
}

for ($i = 0; $i -lt $nodesOrder.Length; $i++) {
    $r = $i + 2
    $name = $nodesOrder[$i]
    $vals = $nodeData[$name]
    $nodesSheet.Range("A$r").Value2 = $name
    $nodesSheet.Range("B$r").Value2 = $vals[0]
    $nodesSheet.Range("C$r").Value2 = $vals[1]
    # has_state ("true"/blank) must stay plain text, not auto-coerced to a
    # boolean - prefix with an apostrophe to force text entry.
    if ($vals[2] -eq "true") {
        $nodesSheet.Range("D$r").Value2 = "'true"
    } else {
        $nodesSheet.Range("D$r").Value2 = $vals[2]
    }
    $nodesSheet.Range("E$r").Value2 = $vals[3]
    $nodesSheet.Range("F$r").Value2 = $vals[4]
    $nodesSheet.Range("G$r").Value2 = $vals[5]
}

# ---------------------------------------------------------------------
# 3) Definition: re-order the node block (rows 7..18) to the same order.
# ---------------------------------------------------------------------
$defSheet = $wb.Worksheets.Item("Definition")

for ($i = 0; $i -lt $nodesOrder.Length; $i++) {
    $r = $i + 7
    $defSheet.Range("A$r").Value2 = $nodesOrder[$i]
    $defSheet.Range("B$r").Value2 = "node"
}

Write-Output "done"
